$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Cell B11 ("R40") is replaced with the text value "1".
# A leading apostrophe forces Excel to store it as text (not a number),
# matching the shared-string / t="s" cell type seen in the target workbook.
$ws.Range("B11").Value = "'1"
